# Updated symbol list (coinranking data refresh) — applies the cell-level
# changes described in the commit diff: price/volume/hour updates for all
# rows, plus a few coins that moved rank (B/C columns) and new price rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'256.81"
$ws.Range("E2").Value = "'4.49%"
$ws.Range("G2").Value = "'4"
$ws.Range("D3").Value = "'27.93"
$ws.Range("E3").Value = "'-5.16%"
$ws.Range("G3").Value = "'4"
$ws.Range("D4").Value = "'5.293"
$ws.Range("E4").Value = "'2.49%"
$ws.Range("G4").Value = "'4"
$ws.Range("D5").Value = "'0.05822"
$ws.Range("E5").Value = "'0.82%"
$ws.Range("G5").Value = "'4"
$ws.Range("D6").Value = "'6.710"
$ws.Range("E6").Value = "'1.59%"
$ws.Range("G6").Value = "'4"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "'0.8720"
$ws.Range("E7").Value = "'1.52%"
$ws.Range("G7").Value = "'4"
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").Value = "'0.9046"
$ws.Range("E8").Value = "'4.94%"
$ws.Range("G8").Value = "'4"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1412"
$ws.Range("E9").Value = "'3.59%"
$ws.Range("G9").Value = "'4"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "'0.07189"
$ws.Range("E10").Value = "'2.34%"
$ws.Range("G10").Value = "'4"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").Value = "'0.03155"
$ws.Range("E11").Value = "'5.52%"
$ws.Range("G11").Value = "'4"
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").Value = "'0.09236"
$ws.Range("E12").Value = "'-1.41%"
$ws.Range("G12").Value = "'4"
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").Value = "'0.001541"
$ws.Range("E13").Value = "'-0.18%"
$ws.Range("G13").Value = "'4"
$ws.Range("B14").Value = "TigerCash"
$ws.Range("C14").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D14").Value = "'0.006035"
$ws.Range("E14").Value = "'1.02%"
$ws.Range("G14").Value = "'4"
$ws.Range("B15").Value = "LEO"
$ws.Range("C15").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D15").Value = "'3.506"
$ws.Range("E15").Value = "'0.56%"
$ws.Range("G15").Value = "'4"
$ws.Range("B16").Value = "GateToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D16").Value = "'3.232"
$ws.Range("E16").Value = "'2.64%"
$ws.Range("G16").Value = "'4"
$ws.Range("B17").Value = "BTSEToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D17").Value = "'2.271"
$ws.Range("E17").Value = "'4.69%"
$ws.Range("G17").Value = "'4"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0006031"
$ws.Range("E18").Value = "'0.88%"
$ws.Range("G18").Value = "'4"
$ws.Range("E19").Value = "'-2.33%"
$ws.Range("G19").Value = "'4"
$ws.Range("D20").Value = "'0.03418"
$ws.Range("E20").Value = "'3.13%"
$ws.Range("G20").Value = "'4"
$ws.Range("D21").Value = "'0.1314"
$ws.Range("E21").Value = "'2.47%"
$ws.Range("G21").Value = "'4"
$ws.Range("D22").Value = "'3.525"
$ws.Range("E22").Value = "'11.15%"
$ws.Range("G22").Value = "'4"
$ws.Range("D23").Value = "'0.04157"
$ws.Range("E23").Value = "'0.20%"
$ws.Range("G23").Value = "'4"
$ws.Range("G24").Value = "'4"
$ws.Range("D25").Value = "'0.001227"
$ws.Range("E25").Value = "'0.12%"
$ws.Range("G25").Value = "'4"
$ws.Range("D26").Value = "'0.004987"
$ws.Range("E26").Value = "'20.60%"
$ws.Range("G26").Value = "'4"
$ws.Range("D27").Value = "'0.0001198"
$ws.Range("E27").Value = "'-0.92%"
$ws.Range("G27").Value = "'4"
$ws.Range("D28").Value = "'0.0001935"
$ws.Range("E28").Value = "'33.60%"
$ws.Range("G28").Value = "'4"
$ws.Range("G29").Value = "'4"
$ws.Range("G30").Value = "'4"
$ws.Range("G31").Value = "'4"
$ws.Range("G32").Value = "'4"
$ws.Range("G33").Value = "'4"
$ws.Range("G34").Value = "'4"
$ws.Range("G35").Value = "'4"
$ws.Range("G36").Value = "'4"
$ws.Range("G37").Value = "'4"
$ws.Range("G38").Value = "'4"
$ws.Range("G39").Value = "'4"
$ws.Range("D40").Value = "'0.03868"
$ws.Range("E40").Value = "'3.58%"
$ws.Range("G40").Value = "'4"
$ws.Range("D41").Value = "'0.005740"
$ws.Range("E41").Value = "'63.70%"
$ws.Range("G41").Value = "'4"
$ws.Range("D42").Value = "'0.1098"
$ws.Range("E42").Value = "'2.52%"
$ws.Range("G42").Value = "'4"
$ws.Range("D43").Value = "'0.002197"
$ws.Range("E43").Value = "'-9.92%"
$ws.Range("G43").Value = "'4"
$ws.Range("D44").Value = "'0.009989"
$ws.Range("E44").Value = "'18.12%"
$ws.Range("G44").Value = "'4"
$ws.Range("D45").Value = "'0.00005282"
$ws.Range("E45").Value = "'-0.24%"
$ws.Range("G45").Value = "'4"
$ws.Range("E46").Value = "'-0.11%"
$ws.Range("G46").Value = "'4"
$ws.Range("D47").Value = "'0.08487"
$ws.Range("E47").Value = "'46.38%"
$ws.Range("G47").Value = "'4"
$ws.Range("D48").Value = "'0.002191"
$ws.Range("E48").Value = "'0.75%"
$ws.Range("G48").Value = "'4"
$ws.Range("E49").Value = "'-0.11%"
$ws.Range("G49").Value = "'4"
$ws.Range("D50").Value = "'0.0001997"
$ws.Range("E50").Value = "'-0.11%"
$ws.Range("G50").Value = "'4"
$ws.Range("G51").Value = "'4"
